$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 7, pushing the existing
# rows 7-94 down to 8-95 (same as the prior week's edit pattern).
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new week's data.
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Vega Monumental Concepción"
$ws.Range("C7").Value = "Bíobío"
$ws.Range("D7").Value = 44552
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 100112043
$ws.Range("G7").Value = "Pepino ensalada"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7500
$ws.Range("N7").Value = "$/caja 60 unidades"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 125
$ws.Range("Q7").Value = 60
$ws.Range("R7").Value = "Hortaliza"
